$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.517.54'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.226.64'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.58'
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.90'
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("D9").Value = '3.223.49'
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("E11").Value = '  -2.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.411'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").Value = '3.779.99'
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.79'
$ws.Range("E15").Value = '  -3.18%  '
$ws.Range("D16").Value = '67.553.64'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("D18").Value = '3.213.48'
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("E19").Value = '  -2.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.44'
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '395.37'
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("E22").Value = '  -2.22%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.17'
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000118'
$ws.Range("E26").Value = '  -2.71%  '
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.56'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("E31").Value = '  -4.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.65'
$ws.Range("E32").Value = '  -1.40%  '
$ws.Range("E33").Value = '  -4.26%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  -2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.37'
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("E37").Value = '  -4.83%  '
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.24'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.801'
$ws.Range("E40").Value = '  -4.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.56'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("E43").Value = '  -5.19%  '
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.51'
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("D46").Value = '2.592.90'
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.53'
$ws.Range("E47").Value = '  -4.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '333.64'
$ws.Range("E49").Value = '  -2.72%  '
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("E51").Value = '  -2.02%  '
